$d = $word.ActiveDocument

# Find.Execute signature:
#   FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#   MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace
# Wrap=1 (wdFindContinue), Replace=2 (wdReplaceAll)

# H1 title (appears again near the bottom, restated in bold) -> new SEO title
$d.Content.Find.Execute("Play Football Cash Collect for Free - Review by Slot Game Writer", $true, $false, $false, $false, $false, $true, 1, $false, "Play Football Cash Collect Free Online", 2)

# "What we like" bullet list
$d.Content.Find.Execute("Four jackpots up for grabs", $true, $false, $false, $false, $false, $true, 1, $false, "Realistic football stadium theme", 2)
$d.Content.Find.Execute("Free spin feature", $true, $false, $false, $false, $false, $true, 1, $false, "Exciting free spin and bonus game features", 2)
$d.Content.Find.Execute("Football stadium theme and design", $true, $false, $false, $false, $false, $true, 1, $false, "High payout potential with the jackpot feature", 2)
$d.Content.Find.Execute("Cash Collect feature adds excitement", $true, $false, $false, $false, $false, $true, 1, $false, "Engaging gameplay with medium to high volatility", 2)

# "What we don't like" bullet list
$d.Content.Find.Execute("Lower than average RTP of 94.91%", $true, $false, $false, $false, $false, $true, 1, $false, "Limited symbol animations", 2)
$d.Content.Find.Execute("No special animations on symbols", $true, $false, $false, $false, $false, $true, 1, $false, "RTP slightly below industry average", 2)

# Meta description (italic) at the very end of the document
$d.Content.Find.Execute("Read our review of Football Cash Collect, a Playtech slot game with 5x3 layout, 30 paylines, free spins, and four jackpots. Play for free today.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Football Cash Collect and play this slot game for free.", 2)
